# ProjectLog update: log entries for guard-fix work, plus a follow-up
# row describing the remaining "humping leg" issue.
#
# Order of operations matters for the on-disk shared-strings table
# (new unique strings are appended in first-use order), so new cells
# are written in the same sequence the author appears to have used:
# row 27 (A-D) first, then the E26 note update, then E27 last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F25: "NPC dialogue added" -> "Text on signs broken"
$ws.Range("F25").Value = "Text on signs broken"

# New row 27: "Guard kind of fixed" entry
$ws.Range("A27").Value = "Guard kind of fixed"

# Dates are stored as literal text (not serial dates) in this log, so
# force the Text number format before assigning, then drop the style
# again so the cell is plain (no explicit style index), matching the
# rest of the sheet's unstyled data cells.
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "02/19/2017"
$ws.Range("B27").ClearFormats()

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "02/20/2017"
$ws.Range("C27").ClearFormats()

$ws.Range("D27").Value = "N/A"

# E26: expand the NPC dialogue note
$ws.Range("E26").Value = "NPC dialogue added, need to disable movement when talking to NPC"

# E27: note about the guard still needing work
$ws.Range("E27").Value = "Need to make guard stop a bit away and stop all movement, as at the moment the guard will be humping leg"

# Leave the selection where the author left it last.
$ws.Range("E27").Select()
